# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed "K" values (column G) keyed by row number.
$gValues = @{
    2  = 2
    4  = 1
    5  = 2
    7  = 0
    8  = 2
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    23 = 0
    24 = 0
    25 = 2
    26 = 1
    27 = 1
    28 = 2
    29 = 0
    30 = 2
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 3
    37 = 1
    39 = 2
    40 = 1
    41 = 0
    42 = 0
    43 = 0
    44 = 2
    45 = 1
    46 = 1
    47 = 2
    48 = 1
    49 = 0
    50 = 1
    51 = 2
    52 = 2
    53 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
